# Update the date heading and the multiplication answers in the table.
# All target strings are unique in the document, but one new value
# ("48×34=1632") happens to coincide with an existing, not-yet-updated
# value elsewhere in the table, so we must perform that swap before the
# replacement that (re)introduces "48×34=1632", to avoid a double hit
# when Find/Replace runs across the whole document.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Heading date
Replace-Text "2025-12-20 Saturday" "2025-12-21 Sunday"

# Row 1
Replace-Text "19×35=665" "90×48=4320"
Replace-Text "87×25=2175" "76×31=2356"
Replace-Text "75×83=6225" "26×26=676"
Replace-Text "39×64=2496" "23×99=2277"
Replace-Text "12×86=1032" "58×91=5278"

# Row 2 (second populated row, "36×85=3060" group)
Replace-Text "36×85=3060" "28×75=2100"
Replace-Text "38×16=608" "74×24=1776"
Replace-Text "59×18=1062" "88×41=3608"
Replace-Text "43×14=602" "92×30=2760"
Replace-Text "55×18=990" "41×57=2337"

# Row 3 ("52×72=3744" group).
# Important: do the "48×34=1632" -> "22×53=1166" swap BEFORE writing a
# new "48×34=1632" value (from "52×72=3744"), otherwise the later
# whole-document Replace-All would also re-match the freshly written cell.
Replace-Text "48×34=1632" "22×53=1166"
Replace-Text "52×72=3744" "48×34=1632"
Replace-Text "82×40=3280" "50×65=3250"
Replace-Text "86×31=2666" "64×70=4480"
Replace-Text "11×34=374" "12×87=1044"

# Row 4 ("69×57=3933" group)
Replace-Text "69×57=3933" "47×58=2726"
Replace-Text "85×57=4845" "63×98=6174"
Replace-Text "99×13=1287" "72×19=1368"
Replace-Text "28×71=1988" "11×13=143"
Replace-Text "72×71=5112" "69×63=4347"

# Row 5 ("75×78=5850" group)
Replace-Text "75×78=5850" "34×93=3162"
Replace-Text "70×69=4830" "64×27=1728"
Replace-Text "54×65=3510" "99×81=8019"
Replace-Text "66×65=4290" "96×73=7008"
Replace-Text "87×92=8004" "14×23=322"
